# "updated parser to take in metadata"
#
# - Cohort sheet ("sheet2"): insert a new "Metadata" / "X" row above the
#   existing "About" row, keep "About" + its long description, rename the
#   template/flags rows and update the template text, resize the Template
#   data row.
# - Active tab moves from "Freshman Profile_SAT Score" to "Cohort".
# - A few sheets get their remembered cell selection moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cohort")

# Insert a new row 20 (pushes the former rows 20-25 down to 21-26, carrying
# their formatting/height down with them since Excel's default insert copies
# the format of the row above).
$ws.Rows(20).Insert()

# New row: Metadata / X
$ws.Range("A20").Value = "Metadata"
$ws.Range("B20").Value = "X"

# Row 21 keeps "About" + the long description untouched (shifted down from 20).

# Row 22 is the old Template row (shifted down from 21): update the template
# wording and resize the row to fit it.
$ws.Range("B22").Value = "The `$aggregation of students {who graduated `$range `$number years and `$range `$number years} in the `$initial_final `$year cohort is [value]"
$ws.Rows(22).RowHeight = 28.8

# Rows 23-26 (old "Operation Allowed?"/"Sum up?"/"Answer Range?"/"Percentage",
# shifted down from 22-25): rename the flag labels.
$ws.Range("A23").Value = "Operation-Allowed?"
$ws.Range("A24").Value = "Sum-Allowed?"
$ws.Range("A25").Value = "Range-Alllowed?"
$ws.Range("A26").Value = "Percentage-Allowed?"

# --- Selections / active sheet -------------------------------------------
# Move the remembered selection on a couple of other sheets first (each of
# these temporarily becomes the active sheet, same as a user clicking
# through tabs)...
$wsBasis = $wb.Worksheets.Item("Basis For Selection")
$wsBasis.Range("A21").Select()

$wsRank = $wb.Worksheets.Item("Freshman Profile_Class Rank GPA")
$wsRank.Range("D25").Select()

# ...then finish on the Cohort sheet so it ends up as the active tab, with
# its new selection at A20 - matching the xlsx's final activeTab/tabSelected.
$ws.Range("A20").Select()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
